$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item("Summary")
$logreg = $wb.Worksheets.Item("LogisticRegression - Obesity")

# --- Insert two new rows at the top of Summary (title row + blank spacer row) ---
$summary.Rows("1:2").Insert()

# New title cell A1: re-use the workbook's existing bold/red "section title" look
# (same look as LogisticRegression - Obesity!A2) by copying its format - this keeps
# the style table clean instead of building it up one font tweak at a time.
$summary.Range("A1").Value = "Using LR, one hot encoding and ngram(1,2)"
$logreg.Range("A2").Copy()
$summary.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column A is now wider to fit the new title
$summary.Columns("A").ColumnWidth = 38.3

# Update the remembered selection on the Summary sheet
$summary.Range("J9").Select()

# --- Add a new blank worksheet ("Sheet1") right after Summary ---
$newSheet = $wb.Worksheets.Add($null, $summary)
$newSheet.Range("G5").Select()

# Make Summary the active/visible tab again (it was tabSelected in the source file)
$summary.Activate()
